$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Basis")
$ws.Activate()

$ws.Range("A7").Value = "Afgerond"
$ws.Range("A8").Value = "Afgerond"
$ws.Range("A9").Value = "Afgerond"
$ws.Range("A11").Value = "Afgerond"
$ws.Range("A12").Value = "Mee bezig"
$ws.Range("A13").Value = "Mee bezig"
$ws.Range("A15").Value = "Afgerond"
$ws.Range("A16").Value = "Afgerond"
$ws.Range("A17").Value = "Afgerond"
$ws.Range("A18").Value = "Afgerond"

$ws.Range("D19").Select()
